$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 ("Save Contract as PDF to Membership Accounting"),
# shifting it (and everything below) down by one row.
$ws.Rows.Item(4).Insert()

# New row 4 holds the new "Colorado contract" line item.
$ws.Range("D4").Value = "Colorado contract"

# Match the author's final selection.
$ws.Range("D4").Select()
